$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update L33:L40 values from 90 to 0
foreach ($r in 33..40) {
    $ws.Cells.Item($r, 12).Value = 0
}

# Update the view: scroll position (topLeftCell) and selection
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L41").Select()
